# Add new columns I ("I0") and J ("IF") to the sheet, mirroring the
# style of the existing header row and populating data rows 2-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - same style as the other header cells (bold + border)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data for columns I (I0) and J (IF), rows 2-53
$data = @(
    @{Row=2; I=8; J=8},
    @{Row=3; I=8; J=8},
    @{Row=4; I=9; J=9},
    @{Row=5; I=6; J=7},
    @{Row=6; I=8; J=8},
    @{Row=7; I=9; J=9},
    @{Row=8; I=9; J=9},
    @{Row=9; I=7; J=7},
    @{Row=10; I=9; J=9},
    @{Row=11; I=8; J=8},
    @{Row=12; I=8; J=8},
    @{Row=13; I=9; J=9},
    @{Row=14; I=7; J=7},
    @{Row=15; I=7; J=7},
    @{Row=16; I=5; J=6},
    @{Row=17; I=7; J=7},
    @{Row=18; I=7; J=7},
    @{Row=19; I=7; J=8},
    @{Row=20; I=6; J=6},
    @{Row=21; I=8; J=8},
    @{Row=22; I=8; J=8},
    @{Row=23; I=10; J=10},
    @{Row=24; I=7; J=7},
    @{Row=25; I=4; J=5},
    @{Row=26; I=6; J=6},
    @{Row=27; I=8; J=8},
    @{Row=28; I=8; J=8},
    @{Row=29; I=7; J=7},
    @{Row=30; I=7; J=7},
    @{Row=31; I=6; J=6},
    @{Row=32; I=7; J=7},
    @{Row=33; I=6; J=6},
    @{Row=34; I=8; J=8},
    @{Row=35; I=6; J=6},
    @{Row=36; I=6; J=6},
    @{Row=37; I=9; J=9},
    @{Row=38; I=8; J=8},
    @{Row=39; I=7; J=7},
    @{Row=40; I=6; J=6},
    @{Row=41; I=8; J=8},
    @{Row=42; I=8; J=8},
    @{Row=43; I=9; J=9},
    @{Row=44; I=9; J=9},
    @{Row=45; I=9; J=9},
    @{Row=46; I=9; J=9},
    @{Row=47; I=6; J=6},
    @{Row=48; I=6; J=6},
    @{Row=49; I=5; J=5},
    @{Row=50; I=4; J=4},
    @{Row=51; I=4; J=4},
    @{Row=52; I=4; J=4},
    @{Row=53; I=4; J=4}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
}

Write-Output "Added columns I0 and J (IF) for rows 1-53"
